# Fruta / hortaliza, semanal
# Weekly refresh of the "Mora" (blackberry) price sheet: the per-row
# observations (date, volume, min/max/weighted price, origin, $/Kg) are
# re-shuffled across rows 2-15 to reflect the new weekly pull. Row 6 is
# unchanged. Only the cells that actually differ from the current values
# are touched below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (now carries what used to be row 12's data)
$ws.Range("D2").Value = 44586
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("S2").Value = 2500

# Row 3 (now carries what used to be row 13's data)
$ws.Range("D3").Value = 44174
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 3200
$ws.Range("O3").Value = 3200
$ws.Range("P3").Value = 3200
$ws.Range("S3").Value = 1600

# Row 4 (now carries what used to be row 10's data)
$ws.Range("D4").Value = 44194
$ws.Range("M4").Value = 120
$ws.Range("R4").Value = "Provincia de Linares"

# Row 5 (now carries what used to be row 4's data)
$ws.Range("D5").Value = 44232
$ws.Range("M5").Value = 200
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 1500

# Row 6 is unchanged.

# Row 7 (now carries what used to be row 2's data)
$ws.Range("D7").Value = 44231
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 3400
$ws.Range("O7").Value = 3400
$ws.Range("P7").Value = 3400
$ws.Range("S7").Value = 1700

# Row 8 (now carries what used to be row 7's data)
$ws.Range("D8").Value = 44237
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 3600
$ws.Range("O8").Value = 4000
$ws.Range("P8").Value = 3800
$ws.Range("S8").Value = 1900

# Row 9 (now carries what used to be row 8's data)
$ws.Range("D9").Value = 44617
$ws.Range("M9").Value = 90
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 6500
$ws.Range("P9").Value = 6500
$ws.Range("S9").Value = 3250

# Row 10 (now carries what used to be row 14's data)
$ws.Range("D10").Value = 44582
$ws.Range("M10").Value = 380
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 5000
$ws.Range("P10").Value = 5000
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 2500

# Row 11 (now carries what used to be row 3's data)
$ws.Range("D11").Value = 44238

# Row 12 (now carries what used to be row 15's data)
$ws.Range("D12").Value = 44168
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("R12").Value = "Provincia de Linares"
$ws.Range("S12").Value = 4000

# Row 13 (now carries what used to be row 5's data)
$ws.Range("D13").Value = 44188
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 3000
$ws.Range("O13").Value = 3400
$ws.Range("P13").Value = 3240
$ws.Range("R13").Value = "Provincia de Linares"
$ws.Range("S13").Value = 1620

# Row 14 (now carries what used to be row 9's data)
$ws.Range("D14").Value = 44533
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 4000
$ws.Range("O14").Value = 4000
$ws.Range("P14").Value = 4000
$ws.Range("S14").Value = 2000

# Row 15 (now carries what used to be row 11's data)
$ws.Range("D15").Value = 44236
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 3600
$ws.Range("O15").Value = 4000
$ws.Range("P15").Value = 3800
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 1900
